$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = "Code"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Definition"
$ws.Range("E1").Value = "Entity_Applied_To"
$ws.Range("F1").Value = "Fields_Excluded_From_View"

$ws.Range("A2").Value = "1.0"
$ws.Range("B2").Value = "01"
$ws.Range("C2").Value = "Application Protection"
$ws.Range("D2").Value = "Information about an application is protected until the resources are confirmed as successful. This includes protecting applications which have had a decision other than a full or partial award."
$ws.Range("E2").Value = "Application`nApplication Decision`nOrganisation (if linked to application)`nGroup (if linked to application)`nPerson (if linked to application)"
$ws.Range("F2").Value = "All application related data (including Application Review and Application Decision)`nAll person data linked to the application`nAll organisation data linked to the application`nAll group data linked to the application"

$ws.Range("A3").Value = "1.0"
$ws.Range("B3").Value = "02"
$ws.Range("C3").Value = "Application Review"
$ws.Range("D3").Value = "Information about the review of an application (including scores) is not available publically"
$ws.Range("E3").Value = "Application Review`nPerson (if linked to application)`nOrganisation (if linked to application)`nGroup (if linked to application)"
$ws.Range("F3").Value = "All Application Review data`nAll person data linked to the application review`nAll organisation data linked to the application review`nAll group data linked to the application review"

$ws.Range("A4").Value = "1.0"
$ws.Range("B4").Value = "03"
$ws.Range("C4").Value = "Commercial Financial"
$ws.Range("D4").Value = "Financial information is competitive until research is completed / released"
$ws.Range("E4").Value = "Application`nAward Granted `nResource Distributed`nAward Received `nResource Received"
$ws.Range("F4").Value = "Application | Minimum Resource Quantity`nApplication | Maximum Resource Quantity`nApplication | Resource Value`nAward Granted | Allocated Resource | Minimum Resource Quantity`nAward Granted | Allocated Resource | Maximum Resource Quantity`nAward Granted | Allocated Resource | Resource Value`nResource Distributed | Resource Actual Quantity`nResource Distributed | Resource Value`nAward Received | Allocated Resource | Minimum Resource Quantity`nAward Received | Allocated Resource | Maximum Resource Quantity`nAward Received | Allocated Resource | Resource Value`nResource Received | Resource Actual Quantity`nResource Received | Resource Value"

$ws.Range("A5").Value = "1.0"
$ws.Range("B5").Value = "04"
$ws.Range("C5").Value = "Topic Identifiable"
$ws.Range("D5").Value = "Topics that some members of the public might be antagonistic towards being researched."
$ws.Range("E5").Value = "Application `nAward Granted `nProject`nAward Received`nOutput"
$ws.Range("F5").Value = "Application | Application Title`nApplication | Outcome Goal`nApplication | Proposal Title`nAward Granted | Award Title`nAward Granted | Award Description`nProject | Project Title`nProject | Project Description`nProject | Keywords`nAward Received | Award Title`nAward Received | Award Description`nOutput | Output Title`nOutput | Output description`nOutput Identifiers"

$ws.Range("A6").Value = "1.0"
$ws.Range("B6").Value = "05"
$ws.Range("C6").Value = "Personal Identifiable"
$ws.Range("D6").Value = "Personally identifiable data is protected."
$ws.Range("E6").Value = "Person"
$ws.Range("F6").Value = "Person | Given Name`nPerson | Other Given Names`nPerson | Family Name`nPerson | Date of Birth`nPerson | Local Person ID`nPerson | Sector Person ID`nPerson | Data Owner ID`nPerson | Prior Local Person ID"

$ws.Range("A7").Value = "1.0"
$ws.Range("B7").Value = "06"
$ws.Range("C7").Value = "Personal Demographic"
$ws.Range("D7").Value = "Information that is used to demographically categorise a person and their work is protected"
$ws.Range("E7").Value = "Person"
$ws.Range("F7").Value = "Person | Gender`nPerson | Other Gender Information`nPerson | Ethnicity`nPerson | Iwi Affiliation"

$ws.Range("A8").Value = "1.0"
$ws.Range("B8").Value = "07"
$ws.Range("C8").Value = "Personal Professional"
$ws.Range("D8").Value = "Information that is used to professionally categorise a person and their work is protected."
$ws.Range("E8").Value = "Person"
$ws.Range("F8").Value = "Person | Career stage`nPerson | Academic Record`nPerson | Recognition`nPerson | Professional Bodies`nPerson | Organisational Affiliation"

$ws.Range("A9").Value = "1.0"
$ws.Range("B9").Value = "08"
$ws.Range("C9").Value = "No Protection Needed"
$ws.Range("D9").Value = "Explicitly stating data is not sensitive"
$ws.Range("E9").Value = "All entities"
$ws.Range("F9").Value = "All fields"

$ws.Range("A10").Value = "1.0"
$ws.Range("B10").Value = "09"
$ws.Range("C10").Value = "Work In Progress Protection"
$ws.Range("D10").Value = "Information that needs protection, providers have some ideas what should be done"
$ws.Range("E10").Value = "All entities"
$ws.Range("F10").Value = "All fields"

$ws.Range("A11").Value = "1.0"
$ws.Range("B11").Value = "10"
$ws.Range("C11").Value = "Unknown protection"
$ws.Range("D11").Value = "Information that needs protection but the specific patterns don't fit."
$ws.Range("E11").Value = "All entities"
$ws.Range("F11").Value = "All fields"
